# "break out stock.yaml completed"
#
# 1. Append the newly-scraped "day" batch (11/11/2024 11:35:10) as rows
#    876-893 on the "day" sheet. bsecode (column D) in this new batch keeps
#    its raw text representation, same as the rest of the freshly scraped
#    row (nsecode, name, timeframe, Date Time are all text too).
# 2. Fix up the "week" sheet rows 529-557 (08/11/2024 batch) whose bsecode
#    (column D) had been stored as text instead of a number - same digits,
#    just retyped as numeric.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Append new rows to the "day" sheet
# ---------------------------------------------------------------------------
$dayWs = $wb.Worksheets.Item("day")

$newRows = @(
    @{ Sr=1;  Code="PAGEIND";    Name="Page Industries Limited";                                  Bse="532827"; Chg=-1.36; Close=47350;              Vol=26937;     Dt="11/11/2024 11:35:10" },
    @{ Sr=2;  Code="APOLLOHOSP"; Name="Apollo Hospitals Enterprise Limited";                       Bse="508869"; Chg=-3.58; Close=7155.45;            Vol=625620;    Dt="11/11/2024 11:35:10" },
    @{ Sr=3;  Code="JKCEMENT";   Name="Jk Cement Limited";                                         Bse="532644"; Chg=-0.42; Close=4064;               Vol=103633;    Dt="11/11/2024 11:35:10" },
    @{ Sr=4;  Code="INDIGO";     Name="Interglobe Aviation Limited";                                Bse="539448"; Chg=0.22;  Close=4011.6;             Vol=231600;    Dt="11/11/2024 11:35:10" },
    @{ Sr=5;  Code="LT";         Name="Larsen & Toubro Limited";                                   Bse="500510"; Chg=-0.86; Close=3628.85;            Vol=926127;    Dt="11/11/2024 11:35:10" },
    @{ Sr=6;  Code="JSWSTEEL";   Name="Jsw Steel Limited";                                         Bse="500228"; Chg=-1.75; Close=979;                Vol=1336153;   Dt="11/11/2024 11:35:10" },
    @{ Sr=7;  Code="TATACONSUM"; Name="TATA Consumer Products Ltd";                                Bse="500800"; Chg=-1.71; Close=975.95;             Vol=1108926;   Dt="11/11/2024 11:35:10" },
    @{ Sr=8;  Code="CONCOR";     Name="Container Corporation Of India Limited";                    Bse="531344"; Chg=0.3;   Close=829.9;              Vol=743364;    Dt="11/11/2024 11:35:10" },
    @{ Sr=9;  Code="SUNTV";      Name="Sun Tv Network Limited";                                    Bse="532733"; Chg=-1.39; Close=744.7;              Vol=474569;    Dt="11/11/2024 11:35:10" },
    @{ Sr=10; Code="GNFC";       Name="Gujarat Narmada Valley Fertilizers And Chemicals Limited";  Bse="500670"; Chg=-1.96; Close=606.55;             Vol=1153728;   Dt="11/11/2024 11:35:10" },
    @{ Sr=11; Code="CHAMBLFERT"; Name="Chambal Fertilizers & Chemicals Limited";                   Bse="500085"; Chg=-0.59; Close=479.8;              Vol=1393659;   Dt="11/11/2024 11:35:10" },
    @{ Sr=12; Code="IGL";        Name="Indraprastha Gas Limited";                                  Bse="532514"; Chg=-0.32; Close=440.95;             Vol=2918340;   Dt="11/11/2024 11:35:10" },
    @{ Sr=13; Code="PETRONET";   Name="Petronet Lng Limited";                                      Bse="532522"; Chg=-0.59; Close=326.95;             Vol=1639447;   Dt="11/11/2024 11:35:10" },
    @{ Sr=14; Code="ONGC";       Name="Oil & Natural Gas Corporation Limited";                     Bse="500312"; Chg=-2.15; Close=256.9;              Vol=8450657;   Dt="11/11/2024 11:35:10" },
    @{ Sr=15; Code="FEDERALBNK"; Name="The Federal Bank  Limited";                                 Bse="500469"; Chg=0.46;  Close=207.73;             Vol=13015971;  Dt="11/11/2024 11:35:10" },
    @{ Sr=16; Code="GMRINFRA";   Name="Gmr Infrastructure Limited";                                Bse="532754"; Chg=-1;    Close=79.48999999999999; Vol=7182172;   Dt="11/11/2024 11:35:10" },
    @{ Sr=17; Code="IDFCFIRSTB"; Name="IDFC First Bank Ltd";                                       Bse="539437"; Chg=1.42;  Close=66.56;              Vol=40877141;  Dt="11/11/2024 11:35:10" },
    @{ Sr=18; Code="IDEA";       Name="Idea Cellular Limited";                                     Bse="532822"; Chg=-0.63; Close=7.83;               Vol=314219334; Dt="11/11/2024 11:35:10" }
)

$startRow = 876
$endRow = $startRow + $newRows.Count - 1

# Pre-format the bsecode column as text for the whole new block so the
# numeric-looking codes ("532827", ...) are stored as strings rather than
# being auto-converted to numbers.
$bseRange = $dayWs.Range("D" + $startRow + ":D" + $endRow)
$bseRange.NumberFormat = "@"

$r = $startRow
foreach ($item in $newRows) {
    $dayWs.Cells.Item($r, 1).Value = $item.Sr
    $dayWs.Cells.Item($r, 2).Value = $item.Code
    $dayWs.Cells.Item($r, 3).Value = $item.Name
    $dayWs.Cells.Item($r, 4).Value = $item.Bse
    $dayWs.Cells.Item($r, 5).Value = $item.Chg
    $dayWs.Cells.Item($r, 6).Value = $item.Close
    $dayWs.Cells.Item($r, 7).Value = $item.Vol
    $dayWs.Cells.Item($r, 8).Value = "day"
    $dayWs.Cells.Item($r, 9).Value = $item.Dt

    $r = $r + 1
}

# Drop the temporary text-number-format style again so the new cells end up
# on the same (unstyled) "Normal" style as every other data row.
$bseRange.Style = "Normal"

# ---------------------------------------------------------------------------
# 2) Fix bsecode (column D) type on the "week" sheet, rows 529-557:
#    was stored as text, should be numeric (same digits, no value change).
# ---------------------------------------------------------------------------
$weekWs = $wb.Worksheets.Item("week")

for ($row = 529; $row -le 557; $row++) {
    $cell = $weekWs.Cells.Item($row, 4)
    $textValue = $cell.Value()
    $cell.Value = [double]$textValue
}
